$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 386.87
$ws.Range("E2").Value = 55.9
$ws.Range("F2").Value = 0.98
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 96
$ws.Range("K2").Value = 66.2
$ws.Range("N2").Value = 54.02451352198364

$ws.Range("F3").Value = 0.13
$ws.Range("K3").Value = 64.59999999999999
$ws.Range("N3").Value = 54.02451352198364

$ws.Range("D4").Value = 4241.1
$ws.Range("E4").Value = 71.7
$ws.Range("F4").Value = 4.43
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 63
$ws.Range("K4").Value = 53.4
$ws.Range("N4").Value = 54.02451352198364
